$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("April")
$ws.Range("A4:A34").NumberFormat = "General"
$wb.Application.Calculate()
Write-Host ("N4=" + $ws.Range("N4").Value2)
Write-Host ("L4=" + $ws.Range("L4").Value2)
Write-Host ("M4=" + $ws.Range("M4").Value2)
Write-Host ("P4=" + $ws.Range("P4").Value2)
Write-Host ("N12=" + $ws.Range("N12").Value2)
Write-Host ("L12=" + $ws.Range("L12").Value2)
